$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data contained in rows 2 and 4, and in rows 3 and 5,
# for the columns that differ between those row pairs (D, K, L, M, N, O, P, Q, S).
# Columns A, B, C, E, F, G, H, I, J, R, T are identical between the swapped
# rows and therefore remain untouched.

# --- Row 2 (becomes what row 4 used to be) ---
$ws.Range("D2").Value = 44902
$ws.Range("K2").Value = "Golden Nugget"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/caja 10 kilos"
$ws.Range("S2").Value = 1500

# --- Row 3 (becomes what row 5 used to be) ---
$ws.Range("D3").Value = 44902
$ws.Range("K3").Value = "Golden Nugget"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("Q3").Value = "$/caja 10 kilos"
$ws.Range("S3").Value = 1300

# --- Row 4 (becomes what row 2 used to be) ---
$ws.Range("D4").Value = 44505
$ws.Range("K4").Value = "Californiana(o)"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1500

# --- Row 5 (becomes what row 3 used to be) ---
$ws.Range("D5").Value = 44505
$ws.Range("K5").Value = "Golden Nugget"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/bandeja 10 kilos"
$ws.Range("S5").Value = 1500
